# Restaurant.xlsx follow-up edit: add "Sales Manager" / "Restaurant Category"
# columns (N, O) to the header row, and populate the single data row with
# "John" / 2 respectively.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Cells.Item(1, 14).Value = "Sales Manager"
$ws.Cells.Item(1, 15).Value = "Restaurant Category"

# New data cells for the existing (only) data row
$ws.Cells.Item(2, 14).Value = "John"
$ws.Cells.Item(2, 15).Value = 2

# Give the new "Restaurant Category" column a sensible, hand-adjusted
# width (close to what Excel auto-applies after typing the content).
$ws.Columns.Item(15).ColumnWidth = 19.6

# Reflect the new active cell / selection left behind by whoever made
# this edit in the Excel UI.
$ws.Range("L14").Select() | Out-Null
